$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 datetime updates (Correspond Handoff/Handback Datetime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 04:01:29"
$wsZhCn.Range("G3").Value = "2016-01-11 04:02:32"

# de-de sheet: row 3 datetime updates (Correspond Handoff/Handback Datetime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 04:01:46"
$wsDeDe.Range("G3").Value = "2016-01-11 04:02:59"
